$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.654.74"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.211.81"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.57"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.01"
$ws.Range("E6").Value = "  +11.21%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.94"
$ws.Range("E10").Value = "  +8.22%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "2.540.81"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.58"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "2.218.68"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "43.613.84"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.96"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.90"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +5.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.57"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.82"
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("E25").Value = "  +21.66%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.53"
$ws.Range("E28").Value = "  +5.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.23"
$ws.Range("E29").Value = "  -9.04%  "
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.07"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0887"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.50"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.33"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  +14.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.23"
$ws.Range("E40").Value = "  -5.59%  "
$ws.Range("E41").Value = "  +7.18%  "
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.21"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.31"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.439"
$ws.Range("E50").Value = "  -5.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.49"
$ws.Range("E51").Value = "  +3.34%  "
